$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains 3 new weekly records at the top of the data block (rows 8-10).
# Insert 3 new blank rows before row 8 -- this shifts the existing rows 8-19
# down to rows 11-22 (dimension grows from A1:R19 to A1:R22).
$ws.Range("A8:R10").Insert()

# --- New row 8: Ají, Americana (o), Primera ---
$ws.Range("A8").Value = 12
$ws.Range("B8").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C8").Value = "Metropolitana"
$ws.Range("D8").Value = 45274
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = 100112021
$ws.Range("G8").Value = "Ají"
$ws.Range("H8").Value = "Americana (o)"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 15
$ws.Range("K8").Value = 50000
$ws.Range("L8").Value = 50000
$ws.Range("M8").Value = 50000
$ws.Range("N8").Value = "`$/caja 25 kilos"
$ws.Range("O8").Value = "Provincia de Huasco"
$ws.Range("P8").Value = 2000
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"

# --- New row 9: Ají, Chilena(o), Primera ---
$ws.Range("A9").Value = 12
$ws.Range("B9").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 45274
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 100112021
$ws.Range("G9").Value = "Ají"
$ws.Range("H9").Value = "Chilena(o)"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 20
$ws.Range("K9").Value = 50000
$ws.Range("L9").Value = 50000
$ws.Range("M9").Value = 50000
$ws.Range("N9").Value = "`$/caja 25 kilos"
$ws.Range("O9").Value = "Provincia de Huasco"
$ws.Range("P9").Value = 2000
$ws.Range("Q9").Value = 25
$ws.Range("R9").Value = "Hortaliza"

# --- New row 10: Ají, Inferno, Primera ---
$ws.Range("A10").Value = 12
$ws.Range("B10").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C10").Value = "Metropolitana"
$ws.Range("D10").Value = 45274
$ws.Range("E10").Value = 13
$ws.Range("F10").Value = 100112021
$ws.Range("G10").Value = "Ají"
$ws.Range("H10").Value = "Inferno"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 10
$ws.Range("K10").Value = 30000
$ws.Range("L10").Value = 30000
$ws.Range("M10").Value = 30000
$ws.Range("N10").Value = "`$/caja 15 kilos"
$ws.Range("O10").Value = "Provincia de Huasco"
$ws.Range("P10").Value = 2000
$ws.Range("Q10").Value = 15
$ws.Range("R10").Value = "Hortaliza"
